$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.2018633540372671
$ws.Range("C2").Value = 0.5652173913043478
$ws.Range("J2").Value = 0.01863354037267081
$ws.Range("P2").Value = 0.1335403726708075
$ws.Range("S2").Value = 0.08074534161490683
$ws.Range("B3").Value = 0.01036269430051814
$ws.Range("C3").Value = 0.03626943005181347
$ws.Range("J3").Value = 0.02072538860103627
$ws.Range("P3").Value = 0.7564766839378239
$ws.Range("S3").Value = 0.1761658031088083
$ws.Range("J4").Value = 0.07894736842105263
$ws.Range("P4").Value = 0.6842105263157895
$ws.Range("S4").Value = 0.2368421052631579
$ws.Range("B6").Value = 0.04938271604938271
$ws.Range("D6").Value = 0.00411522633744856
$ws.Range("E6").Value = 0.00411522633744856
$ws.Range("F6").Value = 0.06995884773662552
$ws.Range("J6").Value = 0.2510288065843622
$ws.Range("O6").Value = 0.0411522633744856
$ws.Range("Q6").Value = 0.1440329218106996
$ws.Range("R6").Value = 0.06584362139917696
$ws.Range("S6").Value = 0.3703703703703703
$ws.Range("B7").Value = 0.1135135135135135
$ws.Range("D7").Value = 0.01081081081081081
$ws.Range("F7").Value = 0.03243243243243243
$ws.Range("J7").Value = 0.1621621621621622
$ws.Range("O7").Value = 0.02162162162162162
$ws.Range("Q7").Value = 0.1891891891891892
$ws.Range("R7").Value = 0.08648648648648649
$ws.Range("S7").Value = 0.3837837837837838
$ws.Range("B8").Value = 0.09049773755656108
$ws.Range("D8").Value = 0.01809954751131222
$ws.Range("F8").Value = 0.05429864253393665
$ws.Range("J8").Value = 0.1176470588235294
$ws.Range("O8").Value = 0.02262443438914027
$ws.Range("Q8").Value = 0.1719457013574661
$ws.Range("R8").Value = 0.1176470588235294
$ws.Range("S8").Value = 0.4072398190045249
$ws.Range("B9").Value = 0.09426229508196721
$ws.Range("D9").Value = 0.02049180327868852
$ws.Range("F9").Value = 0.07377049180327869
$ws.Range("J9").Value = 0.1475409836065574
$ws.Range("O9").Value = 0.03688524590163934
$ws.Range("Q9").Value = 0.1598360655737705
$ws.Range("R9").Value = 0.1024590163934426
$ws.Range("S9").Value = 0.3647540983606558
$ws.Range("B10").Value = 0.1153284671532847
$ws.Range("D10").Value = 0.01897810218978102
$ws.Range("F10").Value = 0.072992700729927
$ws.Range("J10").Value = 0.1175182481751825
$ws.Range("O10").Value = 0.01605839416058394
$ws.Range("Q10").Value = 0.1912408759124088
$ws.Range("R10").Value = 0.0781021897810219
$ws.Range("S10").Value = 0.3897810218978102
$ws.Range("G11").Value = 0.1453287197231834
$ws.Range("J11").Value = 0.08650519031141868
$ws.Range("K11").Value = 0.2041522491349481
$ws.Range("L11").Value = 0.5501730103806228
$ws.Range("S11").Value = 0.01384083044982699
$ws.Range("G12").Value = 0.7701863354037267
$ws.Range("J12").Value = 0.1925465838509317
$ws.Range("K12").Value = 0.0124223602484472
$ws.Range("S12").Value = 0.02484472049689441
$ws.Range("G13").Value = 0.6363636363636364
$ws.Range("J13").Value = 0.2954545454545455
$ws.Range("S13").Value = 0.06818181818181818
$ws.Range("F15").Value = 0.01652892561983471
$ws.Range("H15").Value = 0.1487603305785124
$ws.Range("I15").Value = 0.06198347107438017
$ws.Range("J15").Value = 0.3223140495867768
$ws.Range("K15").Value = 0.0371900826446281
$ws.Range("M15").Value = 0.02066115702479339
$ws.Range("O15").Value = 0.06611570247933884
$ws.Range("S15").Value = 0.3264462809917356
$ws.Range("F16").Value = 0.02392344497607655
$ws.Range("H16").Value = 0.1196172248803828
$ws.Range("I16").Value = 0.1291866028708134
$ws.Range("J16").Value = 0.4162679425837321
$ws.Range("K16").Value = 0.07655502392344497
$ws.Range("M16").Value = 0.01913875598086124
$ws.Range("N16").Value = 0.004784688995215311
$ws.Range("O16").Value = 0.06220095693779904
$ws.Range("S16").Value = 0.1483253588516746
$ws.Range("F17").Value = 0.01138952164009112
$ws.Range("H17").Value = 0.1708428246013667
$ws.Range("I17").Value = 0.1298405466970387
$ws.Range("J17").Value = 0.4396355353075171
$ws.Range("K17").Value = 0.07061503416856492
$ws.Range("M17").Value = 0.01594533029612756
$ws.Range("N17").Value = 0.002277904328018223
$ws.Range("O17").Value = 0.03644646924829157
$ws.Range("S17").Value = 0.1230068337129841
$ws.Range("F18").Value = 0.02325581395348837
$ws.Range("H18").Value = 0.1534883720930233
$ws.Range("I18").Value = 0.1023255813953488
$ws.Range("J18").Value = 0.4186046511627907
$ws.Range("K18").Value = 0.1069767441860465
$ws.Range("M18").Value = 0.0186046511627907
$ws.Range("O18").Value = 0.06046511627906977
$ws.Range("S18").Value = 0.1162790697674419
$ws.Range("F19").Value = 0.02172338884866039
$ws.Range("H19").Value = 0.1976828385228095
$ws.Range("I19").Value = 0.0890658942795076
$ws.Range("J19").Value = 0.3743664011585807
$ws.Range("K19").Value = 0.1064446053584359
$ws.Range("M19").Value = 0.01737871107892831
$ws.Range("O19").Value = 0.07385952208544533
$ws.Range("S19").Value = 0.1194786386676322
